# Git Test.docx edit:
#   "This is the first version of the test document"
#     -> "This is the second version of the test document done on 20/02/23 at 16:00"
#   (split across 4 runs, gramStart/gramEnd proofing marks removed)
#   plus two new empty trailing paragraphs.

$d = $word.ActiveDocument

# The paragraph we need to rewrite is the 3rd paragraph (after the Title and
# the Heading1). It currently reads:
#   "This is the first version of the test " + <proofErr gramStart> + "document" + <proofErr gramEnd>
$oldPara = $d.Paragraphs.Item(3)

# Build the replacement content in a fresh paragraph inserted right after the
# old one, so it inherits the old paragraph's (style-less / Normal) formatting
# rather than the Heading1 formatting of paragraph 2. This also naturally
# avoids carrying over the old proofErr gramStart/gramEnd markers.
$insertionPoint = $oldPara.Range.Duplicate()
$insertionPoint.Collapse(0) | Out-Null
$insertionPoint.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item(4)
$newRange = $newPara.Range.Duplicate()
$newRange.Collapse(0) | Out-Null
$newText = "This is the second version of the test document done on 20/02/23 at 16:00"
$newRange.InsertBefore($newText) | Out-Null

# Split the single inserted run into four runs matching the target:
#   "This is the " | "second" | " version of the test document" | " done on 20/02/23 at 16:00"
# Toggling a character formatting property on a sub-range and then restoring
# it forces the engine to materialize a distinct run there (even though the
# effective formatting ends up identical to its neighbours).
$paraStart = $d.Paragraphs.Item(4).Range.Start

$splitA = $d.Range($paraStart + 12, $paraStart + 18)
$splitA.Bold = 1
$splitA.Bold = 0

$splitB = $d.Range($paraStart + 47, $paraStart + 73)
$splitB.Bold = 1
$splitB.Bold = 0

# Remove the original paragraph (with its stale text and proofErr marks).
$d.Paragraphs.Item(3).Range.Delete() | Out-Null

# Append two new, empty paragraphs after the rewritten one.
$rewritten = $d.Paragraphs.Item(3)
$tail = $rewritten.Range.Duplicate()
$tail.Collapse(0) | Out-Null
$tail.InsertParagraphAfter() | Out-Null

$tail2 = $d.Paragraphs.Item(4).Range.Duplicate()
$tail2.Collapse(0) | Out-Null
$tail2.InsertParagraphAfter() | Out-Null
